$d = $word.ActiveDocument

$replacements = @(
    @{ old = "26×25=650";   new = "72×21=1512" },
    @{ old = "92×98=9016";  new = "62×22=1364" },
    @{ old = "79×36=2844";  new = "74×48=3552" },
    @{ old = "16×76=1216";  new = "33×47=1551" },
    @{ old = "50×93=4650";  new = "65×53=3445" },
    @{ old = "39×66=2574";  new = "79×54=4266" },
    @{ old = "69×88=6072";  new = "21×90=1890" },
    @{ old = "14×25=350";   new = "22×87=1914" },
    @{ old = "56×79=4424";  new = "64×74=4736" },
    @{ old = "53×13=689";   new = "58×98=5684" },
    @{ old = "14×86=1204";  new = "61×39=2379" },
    @{ old = "92×75=6900";  new = "94×58=5452" },
    @{ old = "66×23=1518";  new = "87×28=2436" },
    @{ old = "80×42=3360";  new = "36×46=1656" },
    @{ old = "86×59=5074";  new = "65×30=1950" },
    @{ old = "22×76=1672";  new = "88×60=5280" },
    @{ old = "97×44=4268";  new = "65×57=3705" },
    @{ old = "63×54=3402";  new = "82×81=6642" },
    @{ old = "64×33=2112";  new = "87×47=4089" },
    @{ old = "35×44=1540";  new = "24×61=1464" },
    @{ old = "37×25=925";   new = "82×56=4592" },
    @{ old = "17×71=1207";  new = "85×15=1275" },
    @{ old = "76×37=2812";  new = "33×25=825" },
    @{ old = "11×56=616";   new = "12×36=432" },
    @{ old = "63×99=6237";  new = "89×90=8010" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
